# Updated cryptos list on Thu Jun 27 08:48:56 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to remain stored as text so that
# numeric-looking values (e.g. "569.43", "7.60", "1.00") are not
# reinterpreted by Excel as numbers and reformatted/rounded.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "60.763.64"
$ws.Range("E2").Value = "  -1.25%  "
$ws.Range("D3").Value = "3.371.35"
$ws.Range("E3").Value = "  -0.48%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "569.43"
$ws.Range("D6").Value = "135.74"
$ws.Range("E6").Value = "  -0.72%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "3.369.66"
$ws.Range("E8").Value = "  -0.58%  "
$ws.Range("E9").Value = "  -1.14%  "
$ws.Range("D10").Value = "7.60"
$ws.Range("E11").Value = "  -3.31%  "
$ws.Range("E12").Value = "  -2.95%  "
$ws.Range("D13").Value = "3.945.89"
$ws.Range("E13").Value = "  -0.47%  "
$ws.Range("E14").Value = "  -0.64%  "
$ws.Range("D15").Value = "25.97"
$ws.Range("E15").Value = "  +0.68%  "
$ws.Range("D16").Value = "3.370.13"
$ws.Range("E16").Value = "  -0.54%  "
$ws.Range("E17").Value = "  -3.82%  "
$ws.Range("D18").Value = "60.816.79"
$ws.Range("E18").Value = "  -1.24%  "
$ws.Range("E19").Value = "  -1.19%  "
$ws.Range("D20").Value = "13.74"
$ws.Range("E20").Value = "  -2.95%  "
$ws.Range("D21").Value = "9.19"
$ws.Range("E21").Value = "  -2.25%  "
$ws.Range("D22").Value = "371.54"
$ws.Range("E22").Value = "  -1.59%  "
$ws.Range("D23").Value = "3.507.47"
$ws.Range("E23").Value = "  -0.57%  "
$ws.Range("D24").Value = "0.547"
$ws.Range("E24").Value = "  -1.89%  "
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("D26").Value = "70.66"
$ws.Range("E26").Value = "  -0.85%  "
$ws.Range("D27").Value = "0.0000122"
$ws.Range("E27").Value = "  -2.97%  "
$ws.Range("E28").Value = "  +7.74%  "
$ws.Range("D29").Value = "1.57"
$ws.Range("E29").Value = "  -8.64%  "
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("E31").Value = "  -2.52%  "
$ws.Range("E32").Value = "  -2.76%  "
$ws.Range("E33").Value = "  -2.41%  "
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("D35").Value = "23.26"
$ws.Range("E35").Value = "  -0.79%  "
$ws.Range("E36").Value = "  -4.08%  "
$ws.Range("E37").Value = "  -1.10%  "
$ws.Range("E38").Value = "  -1.13%  "
$ws.Range("D39").Value = "164.40"
$ws.Range("E39").Value = "  -0.75%  "
$ws.Range("D40").Value = "0.0759"
$ws.Range("E41").Value = "  +0.82%  "

# Rows 42, 43 and 45 rotate: the coin that was in row 45 (EnergySwap) moves to
# row 42, the coin that was in row 42 (FirstDigitalUSD) moves to row 43, and
# the coin that was in row 43 (Mantle) moves to row 45. Prices/changes are
# updated to their new values as well.
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").Value = "25.48"
$ws.Range("E42").Value = "  +2.14%  "

$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.00%  "

$ws.Range("E44").Value = "  +0.90%  "

$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").Value = "0.769"
$ws.Range("E45").Value = "  -1.36%  "

$ws.Range("E46").Value = "  -2.23%  "
$ws.Range("D47").Value = "1.16"
$ws.Range("E47").Value = "  -6.53%  "
$ws.Range("D48").Value = "2.512.25"
$ws.Range("E48").Value = "  +7.07%  "
$ws.Range("D49").Value = "23.57"
$ws.Range("E49").Value = "  +4.02%  "
$ws.Range("D50").Value = "6.73"
$ws.Range("E50").Value = "  -1.55%  "
$ws.Range("E51").Value = "  +0.77%  "
